$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("allfund_map")

# Apply Text number format to A3 (new cellXf: numFmtId=49, applyNumberFormat=1 -> s="2")
$ws.Range("A3").NumberFormat = "@"

# Append 85 new fund-manager raw-name -> clean-name mapping rows (rows 157-241)
$ws.Range("A157").Value = "ADVISORS SERIES TRUST"
$ws.Range("B157").Value = "Advisors Series Trust"
$ws.Range("A158").Value = "ADVISORS' INNER CIRCLE FUND"
$ws.Range("B158").Value = "Advisors' Inner Circle"
$ws.Range("A159").Value = "AEGIS FUNDS"
$ws.Range("B159").Value = "Aegis"
$ws.Range("A160").Value = "AFL CIO HOUSING INVESTMENT TRUST"
$ws.Range("B160").Value = "AFL-CIO"
$ws.Range("A161").Value = "ALPS Series Trust"
$ws.Range("B161").Value = "ALPS"
$ws.Range("A162").Value = "AMERICAN PENSION INVESTORS TRUST"
$ws.Range("B162").Value = "American Pension Investors Trust"
$ws.Range("A163").Value = "ARTISAN PARTNERS FUNDS INC"
$ws.Range("B163").Value = "Artisan Partners"
$ws.Range("A164").Value = "Adviser Managed Trust"
$ws.Range("B164").Value = "Adviser Managed Trust"
$ws.Range("A165").Value = "Advisors' Inner Circle Fund II"
$ws.Range("B165").Value = "Advisors' Inner Circle"
$ws.Range("A166").Value = "Advisors' Inner Circle Fund III"
$ws.Range("B166").Value = "Advisors' Inner Circle"
$ws.Range("A167").Value = "Amplify ETF Trust"
$ws.Range("B167").Value = "Amplify ETF"
$ws.Range("A168").Value = "BRIDGEWAY FUNDS INC"
$ws.Range("B168").Value = "Bridgeway"
$ws.Range("A169").Value = "BRUCE FUND INC"
$ws.Range("B169").Value = "Bruce Fund"
$ws.Range("A170").Value = "Barings Funds Trust"
$ws.Range("B170").Value = "Barings"
$ws.Range("A171").Value = "Brighthouse Funds Trust I"
$ws.Range("B171").Value = "Brighthouse"
$ws.Range("A172").Value = "Brighthouse Funds Trust II"
$ws.Range("B172").Value = "Brighthouse"
$ws.Range("A173").Value = "Brinker Capital Destinations Trust"
$ws.Range("B173").Value = "Brinker Capital"
$ws.Range("A174").Value = "Brown Advisory Funds"
$ws.Range("B174").Value = "Brown Advisory"
$ws.Range("A175").Value = "CLEARWATER INVESTMENT TRUST"
$ws.Range("B175").Value = "Clearwater Investment Management"
$ws.Range("A176").Value = "COLORADO BONDSHARES A TAX EXEMPT FUND"
$ws.Range("B176").Value = "Colorado BondShares"
$ws.Range("A177").Value = "COMMONWEALTH INTERNATIONAL SERIES TRUST"
$ws.Range("B177").Value = "Commonwealth Funds"
$ws.Range("A178").Value = "Centre Funds"
$ws.Range("B178").Value = "Centre Funds"
$ws.Range("A179").Value = "Clipper Funds Trust"
$ws.Range("B179").Value = "Clipper"
$ws.Range("A180").Value = "DRIEHAUS MUTUAL FUNDS"
$ws.Range("B180").Value = "Driehaus"
$ws.Range("A181").Value = "EMERGING MARKETS GROWTH FUND INC"
$ws.Range("B181").Value = "Emerging Markets Growth Fund"
$ws.Range("A182").Value = "ETF Series Solutions"
$ws.Range("B182").Value = "ETF Series Solutions"
$ws.Range("A183").Value = "Evermore Funds Trust"
$ws.Range("B183").Value = "Evermore"
$ws.Range("A184").Value = "FIRSTHAND FUNDS"
$ws.Range("B184").Value = "Firsthand Funds"
$ws.Range("A185").Value = "FlexShares Trust"
$ws.Range("B185").Value = "FlexShares"
$ws.Range("A186").Value = "Frost Family of Funds"
$ws.Range("B186").Value = "Frost Investment Advisors"
$ws.Range("A187").Value = "GREAT-WEST FUNDS INC"
$ws.Range("B187").Value = "Great-West"
$ws.Range("A188").Value = "GuideStone Funds"
$ws.Range("B188").Value = "GuideStone"
$ws.Range("A189").Value = "HARBOR FUNDS"
$ws.Range("B189").Value = "Harbor Funds"
$ws.Range("A190").Value = "HC CAPITAL TRUST"
$ws.Range("B190").Value = "HC Capital Trust"
$ws.Range("A191").Value = "HOTCHKIS & WILEY FUNDS /DE/"
$ws.Range("B191").Value = "Hotchkis & Wiley"
$ws.Range("A192").Value = "Jacob Funds Inc."
$ws.Range("B192").Value = "Jacob Funds"
$ws.Range("A193").Value = "KIRR MARBACH PARTNERS FUNDS INC"
$ws.Range("B193").Value = "Kirr, Marbach Partners"
$ws.Range("A194").Value = "LEUTHOLD FUNDS INC"
$ws.Range("B194").Value = "Leuthold"
$ws.Range("A195").Value = "LKCM Funds"
$ws.Range("B195").Value = "LKCM"
$ws.Range("A196").Value = "LoCorr Investment Trust"
$ws.Range("B196").Value = "LoCorr"
$ws.Range("A197").Value = "M FUND INC"
$ws.Range("B197").Value = "Pacific Life"
$ws.Range("A198").Value = "MANNING & NAPIER FUND, INC."
$ws.Range("B198").Value = "Manning & Napier"
$ws.Range("A199").Value = "MATTHEW 25 FUND"
$ws.Range("B199").Value = "Matthew 25"
$ws.Range("A200").Value = "MATTHEWS INTERNATIONAL FUNDS"
$ws.Range("B200").Value = "Matthews"
$ws.Range("A201").Value = "MERGER FUND"
$ws.Range("B201").Value = "Merger Fund"
$ws.Range("A202").Value = "MERGER FUND VL"
$ws.Range("B202").Value = "Merger Fund"
$ws.Range("A203").Value = "MERIDIAN FUND INC"
$ws.Range("B203").Value = "Meridian"
$ws.Range("A204").Value = "MUTUAL OF AMERICA INVESTMENT CORP"
$ws.Range("B204").Value = "Mutual of America"
$ws.Range("A205").Value = "Managed Portfolio Series"
$ws.Range("B205").Value = "Managed Portfolio Series"
$ws.Range("A206").Value = "Meeder Funds"
$ws.Range("B206").Value = "Meeder Funds"
$ws.Range("A207").Value = "Morningstar Funds Trust"
$ws.Range("B207").Value = "Morningstar"
$ws.Range("A208").Value = "Mutual of America Variable Insurance Portfolios, Inc."
$ws.Range("B208").Value = "Mutual of America"
$ws.Range("A209").Value = "NEW ALTERNATIVES FUND"
$ws.Range("B209").Value = "New Alternatives Fund"
$ws.Range("A210").Value = "NORTHEAST INVESTORS TRUST"
$ws.Range("B210").Value = "Northeast Investors Trust"
$ws.Range("A211").Value = "NORTHWESTERN MUTUAL SERIES FUND INC"
$ws.Range("B211").Value = "Northwestern Mutual"
$ws.Range("A212").Value = "North Square Investments Trust"
$ws.Range("B212").Value = "North Square"
$ws.Range("A213").Value = "Northern Lights Fund Trust"
$ws.Range("B213").Value = "Nothern Light"
$ws.Range("A214").Value = "OCM MUTUAL FUND"
$ws.Range("B214").Value = "OCM Funds"
$ws.Range("A215").Value = "OLD WESTBURY FUNDS INC"
$ws.Range("B215").Value = "Old Westbury Funds"
$ws.Range("A216").Value = "PENN SERIES FUNDS INC"
$ws.Range("B216").Value = "Penn Mutual"
$ws.Range("A217").Value = "PPM Funds"
$ws.Range("B217").Value = "PPM Funds"
$ws.Range("A218").Value = "PRIMECAP Odyssey Funds"
$ws.Range("B218").Value = "PRIMECAP"
$ws.Range("A219").Value = "PROFESSIONALLY MANAGED PORTFOLIOS"
$ws.Range("B219").Value = "Osterweis"
$ws.Range("A220").Value = "Pacer Funds Trust"
$ws.Range("B220").Value = "Pacer"
$ws.Range("A221").Value = "QUAKER INVESTMENT TRUST"
$ws.Range("B221").Value = "Quaker Investment Trust"
$ws.Range("A222").Value = "RBB FUND, INC."
$ws.Range("B222").Value = "RBB Fund"
$ws.Range("A223").Value = "RBC FUNDS TRUST"
$ws.Range("B223").Value = "RBC"
$ws.Range("A224").Value = "REYNOLDS FUNDS INC"
$ws.Range("B224").Value = "Reynolds"
$ws.Range("A225").Value = "RMB INVESTORS TRUST"
$ws.Range("B225").Value = "RMB Funds"
$ws.Range("A226").Value = "RiverNorth Funds"
$ws.Range("B226").Value = "RiverNorth"
$ws.Range("A227").Value = "SATURNA INVESTMENT TRUST"
$ws.Range("B227").Value = "Saturna"
$ws.Range("A228").Value = "SPARROW FUNDS"
$ws.Range("B228").Value = "Sparrow"
$ws.Range("A229").Value = "SPIRIT OF AMERICA INVESTMENT FUND INC"
$ws.Range("B229").Value = "Spirit of America"
$ws.Range("A230").Value = "Series Portfolios Trust"
$ws.Range("B230").Value = "Weiss"
$ws.Range("A231").Value = "Symmetry Panoramic Trust"
$ws.Range("B231").Value = "Symmetry Partners"
$ws.Range("A232").Value = "TANAKA FUNDS INC"
$ws.Range("B232").Value = "Tanaka"
$ws.Range("A233").Value = "TIFF INVESTMENT PROGRAM"
$ws.Range("B233").Value = "TIFF"
$ws.Range("A234").Value = "TIMOTHY PLAN"
$ws.Range("B234").Value = "Timothy Plan"
$ws.Range("A235").Value = "TRUST FOR PROFESSIONAL MANAGERS"
$ws.Range("B235").Value = "Trust for Professional Managers"
$ws.Range("A236").Value = "TWEEDY, BROWNE FUND INC."
$ws.Range("B236").Value = "Tweedy, Browne"
$ws.Range("A237").Value = "Tidal ETF Trust"
$ws.Range("B237").Value = "Tidal ETF"
$ws.Range("A238").Value = "Trust for Advised Portfolios"
$ws.Range("B238").Value = "Trust for Advised Portfolios"
$ws.Range("A239").Value = "Vericimetry Funds"
$ws.Range("B239").Value = "Vericimetry"
$ws.Range("A240").Value = "WELLS FARGO MASTER TRUST"
$ws.Range("B240").Value = "Wells Fargo"
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "8"
$ws.Range("ZZ1").Copy()
$ws.Range("A241").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
$ws.Range("B241").Value = "Kinetics Funds"

# Update the _FilterDatabase defined name range to cover the new rows
foreach ($n in $wb.Names) {
    if ($n.Name -eq "allfund_map!_FilterDatabase") {
        $n.RefersTo = "=allfund_map!`$A`$1:`$B`$240"
    }
}

# Move selection to reflect the end of the newly appended data
$ws.Activate()
$ws.Range("B241").Select()
